$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2999
$ws.Range("I137").Value = 2570
$ws.Range("J137").Value = 4000
$ws.Range("K137").Value = 7710
$ws.Range("L137").Value = 12000
$ws.Range("M137").Value = -5160
$ws.Range("N137").Value = -17100

$ws.Range("H138").Value = 2633424
$ws.Range("I138").Value = 1006.80554
$ws.Range("J138").Value = 5002599.5
$ws.Range("K138").Value = 3020.41662
$ws.Range("L138").Value = 15007798.5
$ws.Range("M138").Value = 2119.58338
$ws.Range("N138").Value = -15018078.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 1033.3334
$ws.Range("I41").Value = 1033.3334
$ws.Range("K41").Value = 1033.3334
$ws.Range("M41").Value = -619.3334

$ws.Range("H74").Value = 51546.75
$ws.Range("I74").Value = 72766.14
$ws.Range("J74").Value = 2034.8334
$ws.Range("K74").Value = 72766.14
$ws.Range("L74").Value = 2034.8334
$ws.Range("M74").Value = -71892.14
$ws.Range("N74").Value = -3782.8334

$ws.Range("H77").Value = 51546.75
$ws.Range("I77").Value = 72766.14
$ws.Range("J77").Value = 2034.8334
$ws.Range("K77").Value = 363830.7
$ws.Range("L77").Value = 10174.167
$ws.Range("M77").Value = -359462.7
$ws.Range("N77").Value = -18910.167

$ws.Range("H132").Value = 2562.125
$ws.Range("I132").Value = 2920.111
$ws.Range("J132").Value = 2269.2273
$ws.Range("K132").Value = 8760.332999999999
$ws.Range("L132").Value = 6807.6819
$ws.Range("M132").Value = -6230.332999999999
$ws.Range("N132").Value = -11867.6819

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H44").Value = 20000
$ws.Range("J44").Value = 20000
$ws.Range("L44").Value = 20000
$ws.Range("N44").Value = -20994

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 49999
$ws.Range("J20").Value = 49999
$ws.Range("L20").Value = 49999
$ws.Range("N20").Value = -50471

$ws.Range("H30").Value = 49999
$ws.Range("J30").Value = 49999
$ws.Range("L30").Value = 49999
$ws.Range("N30").Value = -50181

$ws.Range("H56").Value = 8446.5
$ws.Range("I56").Value = 2093
$ws.Range("J56").Value = 14800
$ws.Range("K56").Value = 2093
$ws.Range("L56").Value = 14800
$ws.Range("M56").Value = -1248
$ws.Range("N56").Value = -16490

$ws.Range("H128").Value = 49999
$ws.Range("J128").Value = 49999
$ws.Range("L128").Value = 49999
$ws.Range("N128").Value = -59959

$ws.Range("H132").Value = 4567.9375
$ws.Range("I132").Value = 4882.4
$ws.Range("J132").Value = 4425
$ws.Range("K132").Value = 14647.2
$ws.Range("L132").Value = 13275
$ws.Range("M132").Value = -12117.2
$ws.Range("N132").Value = -18335

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 872.7368
$ws.Range("J5").Value = 1145.909
$ws.Range("L5").Value = 3437.727
$ws.Range("N5").Value = -3661.727

$ws.Range("H13").Value = 373
$ws.Range("I13").Value = 347.6
$ws.Range("J13").Value = 500
$ws.Range("K13").Value = 1042.8
$ws.Range("L13").Value = 1500
$ws.Range("M13").Value = -874.8000000000002
$ws.Range("N13").Value = -1836

$ws.Range("H70").Value = 3831.2273
$ws.Range("I70").Value = 2465.8667
$ws.Range("J70").Value = 6757
$ws.Range("K70").Value = 7397.6001
$ws.Range("L70").Value = 20271
$ws.Range("M70").Value = -7082.6001
$ws.Range("N70").Value = -20901

$ws.Range("H73").Value = 3831.2273
$ws.Range("I73").Value = 2465.8667
$ws.Range("J73").Value = 6757
$ws.Range("K73").Value = 7397.6001
$ws.Range("L73").Value = 20271
$ws.Range("M73").Value = -6305.6001
$ws.Range("N73").Value = -22455

$ws.Range("H82").Value = 3464.8333
$ws.Range("I82").Value = 1750
$ws.Range("J82").Value = 4322.25
$ws.Range("K82").Value = 5250
$ws.Range("L82").Value = 12966.75
$ws.Range("M82").Value = -4844
$ws.Range("N82").Value = -13778.75

$ws.Range("H85").Value = 3464.8333
$ws.Range("I85").Value = 1750
$ws.Range("J85").Value = 4322.25
$ws.Range("K85").Value = 5250
$ws.Range("L85").Value = 12966.75
$ws.Range("M85").Value = -3846
$ws.Range("N85").Value = -15774.75

$ws.Range("H113").Value = 1684042.6
$ws.Range("I113").Value = 6061150
$ws.Range("J113").Value = 539.6923
$ws.Range("K113").Value = 18183450
$ws.Range("L113").Value = 1619.0769
$ws.Range("M113").Value = -18181280
$ws.Range("N113").Value = -5959.0769

$ws.Range("H122").Value = 1160.7391
$ws.Range("I122").Value = 512.5
$ws.Range("J122").Value = 1867.909
$ws.Range("K122").Value = 4612.5
$ws.Range("L122").Value = 16811.181
$ws.Range("M122").Value = -2162.5
$ws.Range("N122").Value = -21711.181

$ws.Range("H132").Value = 2479.25
$ws.Range("I132").Value = 2316.32
$ws.Range("J132").Value = 2849.5454
$ws.Range("K132").Value = 20846.88
$ws.Range("L132").Value = 25645.9086
$ws.Range("M132").Value = -18316.88
$ws.Range("N132").Value = -30705.9086

$ws.Range("H135").Value = 872.7368
$ws.Range("J135").Value = 1145.909
$ws.Range("L135").Value = 10313.181
$ws.Range("N135").Value = -15383.181

$ws.Range("H137").Value = 17814412
$ws.Range("I137").Value = 4240
$ws.Range("J137").Value = 19909726
$ws.Range("K137").Value = 12720
$ws.Range("L137").Value = 59729178
$ws.Range("M137").Value = -7620
$ws.Range("N137").Value = -59739378

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H44").Value = 8000
$ws.Range("J44").Value = 8000
$ws.Range("L44").Value = 8000
$ws.Range("N44").Value = -8912

$ws.Range("H57").Value = 16682
$ws.Range("I57").Value = 10000
$ws.Range("K57").Value = 10000
$ws.Range("M57").Value = -9434

$ws.Range("H127").Value = 34640.92
$ws.Range("J127").Value = 34640.92
$ws.Range("L127").Value = 34640.92
$ws.Range("N127").Value = -44560.92

$ws.Range("H132").Value = 2526.6086
$ws.Range("I132").Value = 2061.7585
$ws.Range("K132").Value = 6185.2755
$ws.Range("M132").Value = -3655.2755

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 8100.5
$ws.Range("J45").Value = 8100.5
$ws.Range("L45").Value = 8100.5
$ws.Range("N45").Value = -9082.5

$ws.Range("H55").Value = 8000.0835
$ws.Range("I55").Value = 3909.6
$ws.Range("J55").Value = 10921.857
$ws.Range("K55").Value = 3909.6
$ws.Range("L55").Value = 10921.857
$ws.Range("M55").Value = -3632.6
$ws.Range("N55").Value = -11475.857

$ws.Range("H61").Value = 14251.4
$ws.Range("I61").Value = 7125
$ws.Range("J61").Value = 19002.334
$ws.Range("K61").Value = 7125
$ws.Range("L61").Value = 19002.334
$ws.Range("M61").Value = -6833
$ws.Range("N61").Value = -19586.334

$ws.Range("H95").Value = 24950
$ws.Range("J95").Value = 24950
$ws.Range("L95").Value = 24950
$ws.Range("N95").Value = -30442
